$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.492.03'
$ws.Range("E2").Value = '  +4.99%  '

# Row 3
$ws.Range("D3").Value = '3.640.82'
$ws.Range("E3").Value = '  +5.06%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.53'
$ws.Range("E5").Value = '  +1.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '194.91'
$ws.Range("E6").Value = '  +3.90%  '

# Row 7
$ws.Range("E7").Value = '  +2.38%  '

# Row 8
$ws.Range("D8").Value = '3.638.73'
$ws.Range("E8").Value = '  +5.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.14%  '

# Row 10
$ws.Range("E10").Value = '  +4.95%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.673'
$ws.Range("E11").Value = '  +3.86%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.53'
$ws.Range("E12").Value = '  +4.42%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000293'
$ws.Range("E13").Value = '  +5.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.99'
$ws.Range("E14").Value = '  +6.12%  '

# Row 15
$ws.Range("D15").Value = '4.219.40'
$ws.Range("E15").Value = '  +4.91%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.96'
$ws.Range("E16").Value = '  +6.23%  '

# Row 17
$ws.Range("D17").Value = '3.638.17'
$ws.Range("E17").Value = '  +4.89%  '

# Row 18
$ws.Range("D18").Value = '70.491.40'
$ws.Range("E18").Value = '  +4.94%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.82'
$ws.Range("E19").Value = '  +5.47%  '

# Row 20
$ws.Range("E20").Value = '  +2.33%  '

# Row 21
$ws.Range("E21").Value = '  +5.20%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.71'
$ws.Range("E22").Value = '  +0.88%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.04'
$ws.Range("E23").Value = '  +12.83%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.36'
$ws.Range("E24").Value = '  +1.45%  '

# Row 25
$ws.Range("E25").Value = '  -0.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.76'
$ws.Range("E26").Value = '  +2.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.18'
$ws.Range("E27").Value = '  +7.68%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.59'
$ws.Range("E28").Value = '  +5.85%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  +6.13%  '

# Row 30
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.96'
$ws.Range("E30").Value = '  +11.05%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.09'
$ws.Range("E31").Value = '  +5.46%  '

# Row 32
$ws.Range("E32").Value = '  +9.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '631.86'
$ws.Range("E33").Value = '  +5.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.31'
$ws.Range("E34").Value = '  +4.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.64'
$ws.Range("E35").Value = '  +2.56%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '40.87'
$ws.Range("E36").Value = '  +11.84%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.416'
$ws.Range("E37").Value = '  +8.04%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0832'
$ws.Range("E38").Value = '  +9.92%  '

# Row 39
$ws.Range("E39").Value = '  -1.66%  '

# Row 40
$ws.Range("E40").Value = '  +0.04%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.60'
$ws.Range("E41").Value = '  +1.56%  '

# Row 42
$ws.Range("D42").Value = '3.317.95'
$ws.Range("E42").Value = '  +2.43%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.92'
$ws.Range("E43").Value = '  +15.75%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.21'
$ws.Range("E44").Value = '  +10.50%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0455'
$ws.Range("E45").Value = '  +5.95%  '

# Row 46
$ws.Range("E46").Value = '  +5.50%  '

# Row 47
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  +0.70%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.139'
$ws.Range("E48").Value = '  +2.64%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.25'

# Row 50
$ws.Range("E50").Value = '  +1.28%  '

# Row 51
$ws.Range("E51").Value = '  -0.05%  '
